# Auto-generated Excel COM-interop script
# Applies market-data refresh values to the Zalera_Profits workbook
# per sheet: currentAveragePrice/LevePrice/LeveProfit columns (H:N)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 8671.333000000001
$ws.Range("I15").Value = 8671.333000000001
$ws.Range("K15").Value = 26013.999
$ws.Range("M15").Value = -25844.999
$ws.Range("H64").Value = 5273
$ws.Range("I64").Value = 5379.7
$ws.Range("K64").Value = 5379.7
$ws.Range("M64").Value = -5131.7
$ws.Range("H67").Value = 5273
$ws.Range("I67").Value = 5379.7
$ws.Range("K67").Value = 5379.7
$ws.Range("M67").Value = -4521.7
$ws.Range("H99").Value = 850
$ws.Range("I99").Value = 300
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = 598
$ws.Range("N99").Value = -10496
$ws.Range("H106").Value = 9000
$ws.Range("I106").Value = 3500
$ws.Range("K106").Value = 3500
$ws.Range("M106").Value = -2869
$ws.Range("H137").Value = 8346182.5
$ws.Range("I137").Value = 19232774
$ws.Range("K137").Value = 57698322
$ws.Range("M137").Value = -57695772

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 24600
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H61").Value = 3737.36
$ws.Range("I61").Value = 3168.762
$ws.Range("J61").Value = 6722.5
$ws.Range("K61").Value = 3168.762
$ws.Range("L61").Value = 6722.5
$ws.Range("M61").Value = -2956.762
$ws.Range("N61").Value = -7146.5
$ws.Range("H74").Value = 14333.333
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 14333.333
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H102").Value = 2398.913
$ws.Range("I102").Value = 2398.913
$ws.Range("K102").Value = 2398.913
$ws.Range("M102").Value = -776.913
$ws.Range("H132").Value = 5586.975
$ws.Range("I132").Value = 3328.56
$ws.Range("J132").Value = 9351
$ws.Range("K132").Value = 9985.68
$ws.Range("L132").Value = 28053
$ws.Range("M132").Value = -7455.68
$ws.Range("N132").Value = -33113
$ws.Range("H136").Value = 3737.36
$ws.Range("I136").Value = 3168.762
$ws.Range("J136").Value = 6722.5
$ws.Range("K136").Value = 9506.286
$ws.Range("L136").Value = 20167.5
$ws.Range("M136").Value = -6956.286
$ws.Range("N136").Value = -25267.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1627.579
$ws.Range("I99").Value = 1656.8889
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1656.8889
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = -158.8888999999999
$ws.Range("N99").Value = -4096
$ws.Range("H107").Value = 3692.2083
$ws.Range("I107").Value = 1980.65
$ws.Range("K107").Value = 1980.65
$ws.Range("M107").Value = -60.65000000000009
$ws.Range("H134").Value = 5829.1665
$ws.Range("I134").Value = 2701.75
$ws.Range("K134").Value = 8105.25
$ws.Range("M134").Value = -5570.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 76927020
$ws.Range("J31").Value = 8249.5
$ws.Range("L31").Value = 8249.5
$ws.Range("N31").Value = -8839.5
$ws.Range("H34").Value = 76927020
$ws.Range("J34").Value = 8249.5
$ws.Range("L34").Value = 8249.5
$ws.Range("N34").Value = -8653.5
$ws.Range("H58").Value = 5657.727
$ws.Range("I58").Value = 4110.625
$ws.Range("K58").Value = 4110.625
$ws.Range("M58").Value = -3907.625
$ws.Range("H86").Value = 6481.6665
$ws.Range("J86").Value = 7060.7144
$ws.Range("L86").Value = 7060.7144
$ws.Range("N86").Value = -9306.714400000001
$ws.Range("H89").Value = 6481.6665
$ws.Range("J89").Value = 7060.7144
$ws.Range("L89").Value = 35303.572
$ws.Range("N89").Value = -46535.572
$ws.Range("H132").Value = 45504.715
$ws.Range("I132").Value = 1697.2273
$ws.Range("J132").Value = 206132.17
$ws.Range("K132").Value = 5091.6819
$ws.Range("L132").Value = 618396.51
$ws.Range("M132").Value = -2561.6819
$ws.Range("N132").Value = -623456.51
$ws.Range("H134").Value = 6849.5835
$ws.Range("I134").Value = 6638.1816
$ws.Range("K134").Value = 19914.5448
$ws.Range("M134").Value = -17379.5448
$ws.Range("H136").Value = 5657.727
$ws.Range("I136").Value = 4110.625
$ws.Range("K136").Value = 12331.875
$ws.Range("M136").Value = -9781.875
$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 50000
$ws.Range("K141").Value = 50000
$ws.Range("M141").Value = -44820

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 460
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 460
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 1380
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -4376
$ws.Range("H131").Value = 17548876
$ws.Range("J131").Value = 8741.5
$ws.Range("L131").Value = 26224.5
$ws.Range("N131").Value = -36304.5
$ws.Range("H132").Value = 1850.4
$ws.Range("I132").Value = 1303.3636
$ws.Range("K132").Value = 11730.2724
$ws.Range("M132").Value = -9200.2724
$ws.Range("H133").Value = 8033
$ws.Range("J133").Value = 8033
$ws.Range("L133").Value = 24099
$ws.Range("N133").Value = -34219

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3807.5
$ws.Range("J80").Value = 4524
$ws.Range("L80").Value = 4524
$ws.Range("N80").Value = -6520
$ws.Range("H83").Value = 3807.5
$ws.Range("J83").Value = 4524
$ws.Range("L83").Value = 22620
$ws.Range("N83").Value = -32604
$ws.Range("H132").Value = 4172.625
$ws.Range("I132").Value = 2304.5
$ws.Range("K132").Value = 6913.5
$ws.Range("M132").Value = -4383.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3033.4
$ws.Range("I16").Value = 1291.375
$ws.Range("K16").Value = 1291.375
$ws.Range("M16").Value = -1121.375
$ws.Range("H46").Value = 6258239
$ws.Range("I46").Value = 20002170
$ws.Range("J46").Value = 10997.454
$ws.Range("K46").Value = 20002170
$ws.Range("L46").Value = 10997.454
$ws.Range("M46").Value = -20001982
$ws.Range("N46").Value = -11373.454
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1298
$ws.Range("N61").Value = -1904
$ws.Range("H100").Value = 10421557
$ws.Range("I100").Value = 13892077
$ws.Range("K100").Value = 13892077
$ws.Range("M100").Value = -13891536
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 9467.308000000001
$ws.Range("I132").Value = 7067.5
$ws.Range("K132").Value = 21202.5
$ws.Range("M132").Value = -18672.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 20123
$ws.Range("I63").Value = 21333.166
$ws.Range("J63").Value = 18912.834
$ws.Range("K63").Value = 21333.166
$ws.Range("L63").Value = 18912.834
$ws.Range("M63").Value = -20709.166
$ws.Range("N63").Value = -20160.834
$ws.Range("H66").Value = 20123
$ws.Range("I66").Value = 21333.166
$ws.Range("J66").Value = 18912.834
$ws.Range("K66").Value = 63999.49800000001
$ws.Range("L66").Value = 56738.50199999999
$ws.Range("M66").Value = -60879.49800000001
$ws.Range("N66").Value = -62978.50199999999
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H96").Value = 11869.909
$ws.Range("I96").Value = 2447
$ws.Range("J96").Value = 23177.4
$ws.Range("K96").Value = 2447
$ws.Range("L96").Value = 23177.4
$ws.Range("M96").Value = -1074
$ws.Range("N96").Value = -25923.4
$ws.Range("H100").Value = 3331.0476
$ws.Range("I100").Value = 3935.6875
$ws.Range("J100").Value = 1396.2
$ws.Range("K100").Value = 7871.375
$ws.Range("L100").Value = 2792.4
$ws.Range("M100").Value = -7330.375
$ws.Range("N100").Value = -3874.4
$ws.Range("H107").Value = 2393.077
$ws.Range("I107").Value = 3496
$ws.Range("K107").Value = 10488
$ws.Range("M107").Value = -8568
$ws.Range("H122").Value = 3746.2942
$ws.Range("I122").Value = 3549.1428
$ws.Range("K122").Value = 10647.4284
$ws.Range("M122").Value = -8197.428400000001
$ws.Range("H126").Value = 6669661
$ws.Range("I126").Value = 7695243.5
$ws.Range("K126").Value = 23085730.5
$ws.Range("M126").Value = -23083260.5
$ws.Range("H132").Value = 5647.8623
$ws.Range("I132").Value = 2666.1875
$ws.Range("K132").Value = 7998.5625
$ws.Range("M132").Value = -5468.5625
$ws.Range("H136").Value = 5077
$ws.Range("I136").Value = 3545.1
$ws.Range("K136").Value = 10635.3
$ws.Range("M136").Value = -8085.299999999999

Write-Host "Applied 238 cell changes across 8 sheets"
